# Auto-generated Excel COM-interop script to apply Halicarnassus_Profits.xlsx data updates
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (59 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 131.23077
$ws.Range("J5").Value = 795
$ws.Range("L5").Value = 795
$ws.Range("N5").Value = -1025
$ws.Range("H17").Value = 3490.5
$ws.Range("I17").Value = 1500
$ws.Range("J17").Value = 3671.4546
$ws.Range("K17").Value = 4500
$ws.Range("L17").Value = 11014.3638
$ws.Range("M17").Value = -4332
$ws.Range("N17").Value = -11350.3638
$ws.Range("H32").Value = 971.25
$ws.Range("I32").Value = 888
$ws.Range("K32").Value = 888
$ws.Range("M32").Value = -562
$ws.Range("H58").Value = 1597.5714
$ws.Range("J58").Value = 2737.5
$ws.Range("L58").Value = 8212.5
$ws.Range("N58").Value = -8512.5
$ws.Range("H62").Value = 9426.5
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 9426.5
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H76").Value = 2061.3333
$ws.Range("J76").Value = 2100
$ws.Range("L76").Value = 2100
$ws.Range("N76").Value = -2730
$ws.Range("H79").Value = 2061.3333
$ws.Range("J79").Value = 2100
$ws.Range("L79").Value = 2100
$ws.Range("N79").Value = -4284
$ws.Range("H86").Value = 2500.5
$ws.Range("I86").Value = 1764.75
$ws.Range("K86").Value = 1764.75
$ws.Range("M86").Value = -641.75
$ws.Range("H89").Value = 2500.5
$ws.Range("I89").Value = 1764.75
$ws.Range("K89").Value = 8823.75
$ws.Range("M89").Value = -3207.75
$ws.Range("H100").Value = 4125.4
$ws.Range("I100").Value = 1202
$ws.Range("K100").Value = 1202
$ws.Range("M100").Value = -661
$ws.Range("H132").Value = 19891.357
$ws.Range("I132").Value = 19039.916
$ws.Range("K132").Value = 57119.74800000001
$ws.Range("M132").Value = -54589.74800000001
$ws.Range("H137").Value = 2333.3333
$ws.Range("I137").Value = 2333.3333
$ws.Range("K137").Value = 6999.999899999999
$ws.Range("M137").Value = -4449.999899999999
$ws.Range("H138").Value = 4182.857
$ws.Range("J138").Value = 4196
$ws.Range("L138").Value = 12588
$ws.Range("N138").Value = -22868

# --- Sheet: ARM (21 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 234
$ws.Range("I5").Value = 167.5
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 167.5
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = -55.5
$ws.Range("N5").Value = -724
$ws.Range("H32").Value = 7728050
$ws.Range("I32").Value = 41725
$ws.Range("J32").Value = 11144194
$ws.Range("K32").Value = 41725
$ws.Range("L32").Value = 11144194
$ws.Range("M32").Value = -41438
$ws.Range("N32").Value = -11144768
$ws.Range("H102").Value = 4255.643
$ws.Range("I102").Value = 2689
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 2689
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -1067
$ws.Range("N102").Value = -13244

# --- Sheet: BSM (18 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 234
$ws.Range("I4").Value = 167.5
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 167.5
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -52.5
$ws.Range("N4").Value = -730
$ws.Range("H22").Value = 400
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -227
$ws.Range("H105").Value = 1699.6
$ws.Range("I105").Value = 1374.5
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 1374.5
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = 372.5
$ws.Range("N105").Value = -6494

# --- Sheet: CRP (50 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7420.737
$ws.Range("I31").Value = 4499.875
$ws.Range("K31").Value = 4499.875
$ws.Range("M31").Value = -4204.875
$ws.Range("H34").Value = 7420.737
$ws.Range("I34").Value = 4499.875
$ws.Range("K34").Value = 4499.875
$ws.Range("M34").Value = -4297.875
$ws.Range("H42").Value = 15028
$ws.Range("I42").Value = 5056
$ws.Range("K42").Value = 5056
$ws.Range("M42").Value = -4463
$ws.Range("H58").Value = 1500
$ws.Range("I58").Value = 1500
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1500
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1297
$ws.Range("N58").ClearContents()
$ws.Range("H88").Value = 15000
$ws.Range("J88").Value = 15000
$ws.Range("L88").Value = 15000
$ws.Range("N88").Value = -15812
$ws.Range("H91").Value = 15000
$ws.Range("J91").Value = 15000
$ws.Range("L91").Value = 15000
$ws.Range("N91").Value = -17808
$ws.Range("H99").Value = 5981.769
$ws.Range("J99").Value = 7322.25
$ws.Range("L99").Value = 7322.25
$ws.Range("N99").Value = -10318.25
$ws.Range("H126").Value = 5981.769
$ws.Range("J126").Value = 7322.25
$ws.Range("L126").Value = 21966.75
$ws.Range("N126").Value = -26906.75
$ws.Range("H134").Value = 3597
$ws.Range("I134").Value = 2956.8572
$ws.Range("J134").Value = 4717.25
$ws.Range("K134").Value = 8870.571599999999
$ws.Range("L134").Value = 14151.75
$ws.Range("M134").Value = -6335.571599999999
$ws.Range("N134").Value = -19221.75
$ws.Range("H136").Value = 1500
$ws.Range("I136").Value = 1500
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4500
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1950
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 16036.923

# --- Sheet: CUL (43 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 596.63635
$ws.Range("I60").Value = 269.5263
$ws.Range("J60").Value = 2668.3333
$ws.Range("K60").Value = 808.5789
$ws.Range("L60").Value = 8004.999899999999
$ws.Range("M60").Value = -557.5789
$ws.Range("N60").Value = -8506.999899999999
$ws.Range("H109").Value = 221385.2
$ws.Range("I109").Value = 275606.5
$ws.Range("J109").Value = 4500
$ws.Range("K109").Value = 826819.5
$ws.Range("L109").Value = 13500
$ws.Range("M109").Value = -825779.5
$ws.Range("N109").Value = -15580
$ws.Range("H115").Value = 3727.8125
$ws.Range("I115").Value = 2440.4546
$ws.Range("K115").Value = 7321.3638
$ws.Range("M115").Value = -6146.3638
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()
$ws.Range("H131").Value = 1364.6666
$ws.Range("I131").Value = 2200
$ws.Range("J131").Value = 1197.6
$ws.Range("K131").Value = 6600
$ws.Range("L131").Value = 3592.8
$ws.Range("M131").Value = -1560
$ws.Range("N131").Value = -13672.8
$ws.Range("H139").Value = 4666.3335
$ws.Range("I139").Value = 2899.5
$ws.Range("J139").Value = 8200
$ws.Range("K139").Value = 8698.5
$ws.Range("L139").Value = 24600
$ws.Range("M139").Value = -3558.5
$ws.Range("N139").Value = -34880
$ws.Range("H140").Value = 2590.2307
$ws.Range("I140").Value = 1837.3
$ws.Range("J140").Value = 5100
$ws.Range("K140").Value = 5511.9
$ws.Range("L140").Value = 15300
$ws.Range("M140").Value = -331.8999999999996
$ws.Range("N140").Value = -25660

# --- Sheet: GSM (23 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4206
$ws.Range("I70").Value = 3999.75
$ws.Range("K70").Value = 3999.75
$ws.Range("M70").Value = -3729.75
$ws.Range("H73").Value = 4206
$ws.Range("I73").Value = 3999.75
$ws.Range("K73").Value = 3999.75
$ws.Range("M73").Value = -3063.75
$ws.Range("H93").Value = 57499.5
$ws.Range("J93").Value = 57499.5
$ws.Range("L93").Value = 57499.5
$ws.Range("N93").Value = -61243.5
$ws.Range("H97").Value = 1481.2
$ws.Range("I97").Value = 1899.3334
$ws.Range("J97").Value = 1302
$ws.Range("K97").Value = 1899.3334
$ws.Range("L97").Value = 1302
$ws.Range("M97").Value = -1403.3334
$ws.Range("N97").Value = -2294
$ws.Range("H108").Value = 84999
$ws.Range("J108").Value = 84999
$ws.Range("L108").Value = 84999
$ws.Range("N108").Value = -92679

# --- Sheet: LTW (25 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4138.875
$ws.Range("I7").Value = 2776.8333
$ws.Range("J7").Value = 8225
$ws.Range("K7").Value = 2776.8333
$ws.Range("L7").Value = 8225
$ws.Range("M7").Value = -2664.8333
$ws.Range("N7").Value = -8449
$ws.Range("H40").Value = 6868.1
$ws.Range("I40").Value = 5526.143
$ws.Range("K40").Value = 5526.143
$ws.Range("M40").Value = -5390.143
$ws.Range("H126").Value = 4138.875
$ws.Range("I126").Value = 2776.8333
$ws.Range("J126").Value = 8225
$ws.Range("K126").Value = 8330.499899999999
$ws.Range("L126").Value = 24675
$ws.Range("M126").Value = -5860.499899999999
$ws.Range("N126").Value = -29615
$ws.Range("H132").Value = 9757.941000000001
$ws.Range("I132").Value = 8574.5
$ws.Range("J132").Value = 12598.2
$ws.Range("K132").Value = 25723.5
$ws.Range("L132").Value = 37794.60000000001
$ws.Range("M132").Value = -23193.5
$ws.Range("N132").Value = -42854.60000000001

# --- Sheet: WVR (7 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 7937.1055
$ws.Range("I136").Value = 7202.7856
$ws.Range("J136").Value = 9993.200000000001
$ws.Range("K136").Value = 21608.3568
$ws.Range("L136").Value = 29979.6
$ws.Range("M136").Value = -19058.3568
$ws.Range("N136").Value = -35079.60000000001
